# Add a new Item row ("Equip_Weapon_1" / 开山斧) to the Item.xlsx "表1" XML table.
#
# This reproduces (as far as the exposed Excel object model allows) the
# change recorded in the commit: a new data row is appended below the
# existing 8 rows, the worksheet/table ranges grow from K8 to K9, and the
# active selection moves to K13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the new row of data (row 9) -----------------------------------
# Columns: ID | ItemType | ItemSubType | Level | ShowName | Desc | Icon |
#          CoolDownTime | OverlayCount | BuyPrice | SalePrice
$ws.Range("A9").Value = "Equip_Weapon_1"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = "开山斧"
$ws.Range("F9").Value = "开山斧武器"

# Icon is stored as text (matches the existing "1017"/"1018"/"1019" style
# used by the other rows), so force a text number format before writing the
# numeric-looking string, otherwise it would be stored as a number.
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "50004"

$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 10000
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 100

# --- Grow the XML table ("表1") so it covers the new row -------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:K9"))

# --- Move the active selection to match the saved view --------------------
$ws.Range("K13").Select()
